$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the record for CARLO PAOLO CANTILLO QUINTANA (employee 1128054869) entirely.
# Deleting the row shifts the remaining worker rows (and the signature block below)
# up by one, which matches the new layout.
$ws.Rows(16).Delete()

# Update the summary figures at the top of the statement.
$ws.Range("E11").Value = 103758
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3

# Rewrite the worker detail rows (now rows 16-18) with the refreshed data:
# JOINER CORTECERO MONTERROZA appears twice (period 2306 then 2307) followed by
# LUIS EDUARDO RICARDO MONTERROSA (period 2312, with an updated "Salario Basico").
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143386223"
$ws.Range("D16").Value = "JOINER CORTECERO MONTERROZA"
$ws.Range("E16").Value = "2306"
$ws.Range("F16").Value = 4640
$ws.Range("G16").Value = 1160000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143386223"
$ws.Range("D17").Value = "JOINER CORTECERO MONTERROZA"
$ws.Range("E17").Value = "2307"
$ws.Range("F17").Value = 43307
$ws.Range("G17").Value = 1160000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047456167"
$ws.Range("D18").Value = "LUIS EDUARDO RICARDO MONTERROSA"
$ws.Range("E18").Value = "2312"
$ws.Range("F18").Value = 55811
$ws.Range("G18").Value = 1528390
